# Translate the "classes" lookup table (column B = human-readable name) from
# Spanish into English, then update the sheet view / selection / column width
# on the "classes" sheet to match the edited workbook, and move the active
# tab from "cameras" to "classes".

$wb = $excel.ActiveWorkbook

$ws_classes = $wb.Worksheets.Item("classes")

# --- Translate the class names (column B, rows 2-17), Spanish -> English ---
# Row order matches the existing "tag" column (A) top to bottom.
$translations = @(
    "Car Front Side",
    "Car Back Side",
    "Car Right Side",
    "Car Left Side",
    "Truck Front Side",
    "Truck Back Side",
    "Truck Right Side",
    "Truck Left Side",
    "Motorcycle Front Side",
    "Motorcycle Back Side",
    "Motorcycle Right Side",
    "Motorcycle Left Side",
    "Bus Front Side",
    "Bus Back Side",
    "Bus Right Side",
    "Bus Left Side"
)

for ($i = 0; $i -lt $translations.Length; $i++) {
    $row = $i + 2
    $ws_classes.Range("B$row").Value = $translations[$i]
}

# --- Widen column B on the "classes" sheet ----------------------------------
# NOTE: the host engine stores/export column widths quantized to 1/6-character
# increments with a fixed +5/6 offset between the COM "ColumnWidth" value and
# the value persisted in the OOXML <col width="..."> attribute. To land as
# close as possible on the target width of 29.42578125 (which exported as
# 29.5, the nearest reachable increment) we request 28.666666666666668.
$ws_classes.Columns.Item(2).ColumnWidth = 28.666666666666668

# --- Move the active tab from "cameras" to "classes" and update selection --
$ws_classes.Activate()
$ws_classes.Range("H17").Select()
